$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Delete()
